# Auto-generated edit script: updates cached computed values in each profession
# sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) as refreshed by the scheduled market-price runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H28").Value = 337.875
$ws.Range("I28").Value = 337.875
$ws.Range("K28").Value = 337.875
$ws.Range("M28").Value = 147.125

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H51").Value = 7925.091
$ws.Range("I51").Value = 6447.8335
$ws.Range("K51").Value = 6447.8335
$ws.Range("M51").Value = -5963.8335

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H74").Value = 6416.048
$ws.Range("I74").Value = 6949.1333
$ws.Range("J74").Value = 5083.3335
$ws.Range("K74").Value = 6949.1333
$ws.Range("L74").Value = 5083.3335
$ws.Range("M74").Value = -6013.1333
$ws.Range("N74").Value = -6955.3335

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H77").Value = 6416.048
$ws.Range("I77").Value = 6949.1333
$ws.Range("J77").Value = 5083.3335
$ws.Range("K77").Value = 34745.6665
$ws.Range("L77").Value = 25416.6675
$ws.Range("M77").Value = -30065.6665
$ws.Range("N77").Value = -34776.6675

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H116").Value = 10830.929
$ws.Range("I116").Value = 7999.4
$ws.Range("K116").Value = 7999.4
$ws.Range("M116").Value = -4557.4

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H138").Value = 2952.5
$ws.Range("J138").Value = 5545.4546
$ws.Range("L138").Value = 16636.3638
$ws.Range("N138").Value = -26916.3638

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 4585.0713
$ws.Range("I2").Value = 3195.4
$ws.Range("K2").Value = 3195.4
$ws.Range("M2").Value = -3082.4

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 3108.9111
$ws.Range("I32").Value = 1701.1621
$ws.Range("J32").Value = 9619.75
$ws.Range("K32").Value = 1701.1621
$ws.Range("L32").Value = 9619.75
$ws.Range("M32").Value = -1414.1621
$ws.Range("N32").Value = -10193.75

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H74").Value = 3094.111
$ws.Range("I74").Value = 1745.2727
$ws.Range("K74").Value = 1745.2727
$ws.Range("M74").Value = -871.2727

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H77").Value = 3094.111
$ws.Range("I77").Value = 1745.2727
$ws.Range("K77").Value = 8726.363499999999
$ws.Range("M77").Value = -4358.363499999999

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H110").Value = 992
$ws.Range("I110").Value = 1025
$ws.Range("J110").Value = 794
$ws.Range("K110").Value = 1025
$ws.Range("L110").Value = 794
$ws.Range("M110").Value = 1020
$ws.Range("N110").Value = -4884

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H116").Value = 4585.0713
$ws.Range("I116").Value = 3195.4
$ws.Range("K116").Value = 3195.4
$ws.Range("M116").Value = -901.4000000000001

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H122").Value = 4401.5386
$ws.Range("J122").Value = 4258.6665
$ws.Range("L122").Value = 12775.9995
$ws.Range("N122").Value = -17675.9995

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H132").Value = 4943.0566
$ws.Range("I132").Value = 5616
$ws.Range("J132").Value = 1653.1111
$ws.Range("K132").Value = 16848
$ws.Range("L132").Value = 4959.3333
$ws.Range("M132").Value = -14318
$ws.Range("N132").Value = -10019.3333

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 4585.0713
$ws.Range("I3").Value = 3195.4
$ws.Range("K3").Value = 3195.4
$ws.Range("M3").Value = -3081.4

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H134").Value = 5143.1665
$ws.Range("I134").Value = 3359.5386
$ws.Range("K134").Value = 10078.6158
$ws.Range("M134").Value = -7543.6158

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 3233.9473
$ws.Range("J31").Value = 4286.727
$ws.Range("L31").Value = 4286.727
$ws.Range("N31").Value = -4876.727

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H34").Value = 3233.9473
$ws.Range("J34").Value = 4286.727
$ws.Range("L34").Value = 4286.727
$ws.Range("N34").Value = -4690.727

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H134").Value = 18132.445
$ws.Range("I134").Value = 16665.666
$ws.Range("K134").Value = 49996.99800000001
$ws.Range("M134").Value = -47461.99800000001

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H7").Value = 115.4
$ws.Range("I7").Value = 80
$ws.Range("K7").Value = 240
$ws.Range("M7").Value = -128

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H50").Value = 8942.904
$ws.Range("I50").Value = 4334.1665
$ws.Range("J50").Value = 9544.044
$ws.Range("K50").Value = 13002.4995
$ws.Range("L50").Value = 28632.132
$ws.Range("M50").Value = -12521.4995
$ws.Range("N50").Value = -29594.132

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H53").Value = 8942.904
$ws.Range("I53").Value = 4334.1665
$ws.Range("J53").Value = 9544.044
$ws.Range("K53").Value = 13002.4995
$ws.Range("L53").Value = 28632.132
$ws.Range("M53").Value = -12521.4995
$ws.Range("N53").Value = -29594.132

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H109").Value = 561.125
$ws.Range("I109").Value = 561.125
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1683.375
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -643.375
$ws.Range("N109").ClearContents()

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H129").Value = 1575.9
$ws.Range("J129").Value = 1833
$ws.Range("L129").Value = 5499
$ws.Range("N129").Value = -15499

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H131").Value = 25806.756
$ws.Range("J131").Value = 1862.037
$ws.Range("L131").Value = 5586.111
$ws.Range("N131").Value = -15666.111

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H5").Value = 15000
$ws.Range("I5").Value = 15000
$ws.Range("K5").Value = 15000
$ws.Range("M5").Value = -14888

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H27").Value = 14997
$ws.Range("J27").Value = 14997
$ws.Range("L27").Value = 14997
$ws.Range("N27").Value = -15329

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H113").Value = 2990.1428
$ws.Range("I113").Value = 2990.1428
$ws.Range("K113").Value = 2990.1428
$ws.Range("M113").Value = -820.1428000000001

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H122").Value = 4457.75
$ws.Range("I122").Value = 4443
$ws.Range("K122").Value = 13329
$ws.Range("M122").Value = -10879

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H82").Value = 1248.7222
$ws.Range("I82").Value = 1374.9231
$ws.Range("K82").Value = 1374.9231
$ws.Range("M82").Value = -1013.9231

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H85").Value = 1248.7222
$ws.Range("I85").Value = 1374.9231
$ws.Range("K85").Value = 1374.9231
$ws.Range("M85").Value = -126.9231

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H122").Value = 3657.8
$ws.Range("I122").Value = 3022.25
$ws.Range("K122").Value = 9066.75
$ws.Range("M122").Value = -6616.75

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H136").Value = 13890482
$ws.Range("I136").Value = 15874680
$ws.Range("K136").Value = 47624040
$ws.Range("M136").Value = -47621490

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H9").Value = 7999
$ws.Range("I9").Value = 7999
$ws.Range("K9").Value = 7999
$ws.Range("M9").Value = -7859

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H41").Value = 21170
$ws.Range("I41").Value = 21170
$ws.Range("K41").Value = 21170
$ws.Range("M41").Value = -20780

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H113").Value = 740
$ws.Range("I113").Value = 549.25
$ws.Range("J113").Value = 994.3333
$ws.Range("K113").Value = 1647.75
$ws.Range("L113").Value = 2982.9999
$ws.Range("M113").Value = 522.25
$ws.Range("N113").Value = -7322.9999

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H136").Value = 1076.5264
$ws.Range("I136").Value = 615
$ws.Range("K136").Value = 1845
$ws.Range("M136").Value = 705
